$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values, forcing text storage to preserve exact formatting
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.495.97'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.640.41'
$ws.Range("D3").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.50'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3789'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '51.89'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3623'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08185'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.237'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9997'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.57'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.467'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.376'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001240'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.637.02'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.24'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06943'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.579'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.55'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.528.46'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.509'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.070'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.22'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '152.33'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.259'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.46'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.812.26'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.108'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.609'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.146'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.49'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02766'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2504'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.08777'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.029'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.07089'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7066'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.353'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.30'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.60'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6554'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9993'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.285'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.968'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07984'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '128.86'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.196'
$ws.Range("D51").Style = "Normal"

# Update Volume(1h) percentages (column E), and Coin/Link (columns B/C) for the
# two rows whose coin order was swapped (rows 39-40: Hedera <-> InternetComputer)
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("E3").Value = '  +0.00%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("E5").Value = '  -0.27%  '
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("E7").Value = '  +0.35%  '
$ws.Range("E8").Value = '  -0.88%  '
$ws.Range("E9").Value = '  -0.53%  '
$ws.Range("E10").Value = '  +0.92%  '
$ws.Range("E11").Value = '  -0.77%  '
$ws.Range("E12").Value = '  -0.49%  '
$ws.Range("E13").Value = '  -1.39%  '
$ws.Range("E14").Value = '  -2.65%  '
$ws.Range("E15").Value = '  +1.11%  '
$ws.Range("E16").Value = '  -1.13%  '
$ws.Range("E17").Value = '  -0.01%  '
$ws.Range("E18").Value = '  +1.13%  '
$ws.Range("E19").Value = '  +0.04%  '
$ws.Range("E20").Value = '  +0.47%  '
$ws.Range("E21").Value = '  -3.33%  '
$ws.Range("E22").Value = '  -0.33%  '
$ws.Range("E23").Value = '  -2.58%  '
$ws.Range("E24").Value = '  +0.31%  '
$ws.Range("E26").Value = '  -4.70%  '
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("E28").Value = '  +0.92%  '
$ws.Range("E29").Value = '  -0.98%  '
$ws.Range("E30").Value = '  -1.83%  '
$ws.Range("E31").Value = '  -0.34%  '
$ws.Range("E32").Value = '  +14.90%  '
$ws.Range("E33").Value = '  -4.02%  '
$ws.Range("E34").Value = '  -7.37%  '
$ws.Range("E35").Value = '  +4.70%  '
$ws.Range("E36").Value = '  -3.19%  '
$ws.Range("E37").Value = '  -1.99%  '
$ws.Range("E38").Value = '  -1.03%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("E39").Value = '  -3.82%  '
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("E40").Value = '  -2.71%  '
$ws.Range("E41").Value = '  -0.84%  '
$ws.Range("E42").Value = '  -1.51%  '
$ws.Range("E43").Value = '  -2.28%  '
$ws.Range("E44").Value = '  -5.16%  '
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("E46").Value = '  -0.31%  '
$ws.Range("E47").Value = '  -3.03%  '
$ws.Range("E49").Value = '  -0.14%  '
$ws.Range("E50").Value = '  +0.65%  '
$ws.Range("E51").Value = '  -1.82%  '
